# Apply the updated cryptocurrency price/volume snapshot to sheet1.
# Generated from the canonical OOXML diff: row-level coin price/percentage
# refreshes plus two rank swaps (rows 38/39 and 50/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range('D2').Value = '63.790.93'
$ws.Range('E2').Value = '  -1.04%  '

# Row 3: update D3, E3
$ws.Range('D3').Value = '2.641.00'
$ws.Range('E3').Value = '  +0.34%  '

# Row 4: update E4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5: update D5, E5
$ws.Range('D5').Value = '580.27'
$ws.Range('E5').Value = '  +0.52%  '

# Row 6: update D6, E6
$ws.Range('D6').Value = '155.65'
$ws.Range('E6').Value = '  -0.33%  '

# Row 7: update E7
$ws.Range('E7').Value = '  +0.07%  '

# Row 8: update D8, E8
$ws.Range('D8').Value = '0.619'
$ws.Range('E8').Value = '  -3.92%  '

# Row 9: update D9, E9
$ws.Range('D9').Value = '2.638.58'
$ws.Range('E9').Value = '  +0.37%  '

# Row 10: update E10
$ws.Range('E10').Value = '  -3.58%  '

# Row 11: update E11
$ws.Range('E11').Value = '  +0.38%  '

# Row 12: update E12
$ws.Range('E12').Value = '  -1.89%  '

# Row 13: update D13, E13
$ws.Range('D13').Value = '0.157'
$ws.Range('E13').Value = '  +1.00%  '

# Row 14: update D14, E14
$ws.Range('D14').Value = '28.44'
$ws.Range('E14').Value = '  -0.28%  '

# Row 15: update D15, E15
$ws.Range('D15').Value = '3.117.08'
$ws.Range('E15').Value = '  +0.43%  '

# Row 16: update D16, E16
$ws.Range('D16').Value = '0.0000182'
$ws.Range('E16').Value = '  -1.81%  '

# Row 17: update D17, E17
$ws.Range('D17').Value = '63.774.66'
$ws.Range('E17').Value = '  -0.79%  '

# Row 18: update D18, E18
$ws.Range('D18').Value = '2.638.84'
$ws.Range('E18').Value = '  +0.52%  '

# Row 19: update D19, E19
$ws.Range('D19').Value = '12.13'
$ws.Range('E19').Value = '  -0.91%  '

# Row 20: update D20, E20
$ws.Range('D20').Value = '7.64'
$ws.Range('E20').Value = '  +3.27%  '

# Row 21: update D21, E21
$ws.Range('D21').Value = '4.52'
$ws.Range('E21').Value = '  -3.30%  '

# Row 22: update D22, E22
$ws.Range('D22').Value = '344.30'
$ws.Range('E22').Value = '  -0.59%  '

# Row 23: update E23
$ws.Range('E23').Value = '  +0.32%  '

# Row 24: update D24, E24
$ws.Range('D24').Value = '68.09'
$ws.Range('E24').Value = '  +0.57%  '

# Row 25: update D25, E25
$ws.Range('D25').Value = '1.88'
$ws.Range('E25').Value = '  +8.05%  '

# Row 26: update D26, E26
$ws.Range('D26').Value = '0.0000109'
$ws.Range('E26').Value = '  -3.18%  '

# Row 27: update D27, E27
$ws.Range('D27').Value = '604.72'
$ws.Range('E27').Value = '  +8.55%  '

# Row 28: update D28, E28
$ws.Range('D28').Value = '9.26'
$ws.Range('E28').Value = '  -1.11%  '

# Row 29: update E29
$ws.Range('E29').Value = '  +2.47%  '

# Row 30: update D30, E30
$ws.Range('D30').Value = '8.09'
$ws.Range('E30').Value = '  +2.44%  '

# Row 31: update E31
$ws.Range('E31').Value = '  -0.19%  '

# Row 32: update E32
$ws.Range('E32').Value = '  -0.11%  '

# Row 33: update E33
$ws.Range('E33').Value = '  -0.06%  '

# Row 34: update D34, E34
$ws.Range('D34').Value = '1.74'
$ws.Range('E34').Value = '  +1.31%  '

# Row 35: update D35, E35
$ws.Range('D35').Value = '6.58'
$ws.Range('E35').Value = '  -0.28%  '

# Row 36: update D36, E36
$ws.Range('D36').Value = '5.42'
$ws.Range('E36').Value = '  +2.60%  '

# Row 37: update D37, E37
$ws.Range('D37').Value = '0.403'
$ws.Range('E37').Value = '  -2.03%  '

# Row 38: update B38, C38, D38, E38
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  +0.08%  '

# Row 39: update B39, C39, D39, E39
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').Value = '19.73'
$ws.Range('E39').Value = '  -1.51%  '

# Row 40: update E40
$ws.Range('E40').Value = '  -1.50%  '

# Row 41: update D41, E41
$ws.Range('D41').Value = '150.87'
$ws.Range('E41').Value = '  -2.41%  '

# Row 42: update D42, E42
$ws.Range('D42').Value = '2.56'
$ws.Range('E42').Value = '  +5.00%  '

# Row 43: update E43
$ws.Range('E43').Value = '  -0.04%  '

# Row 44: update D44, E44
$ws.Range('D44').Value = '41.89'
$ws.Range('E44').Value = '  -0.57%  '

# Row 45: update D45
$ws.Range('D45').Value = '160.55'

# Row 46: update D46, E46
$ws.Range('D46').Value = '24.33'
$ws.Range('E46').Value = '  +6.61%  '

# Row 47: update D47, E47
$ws.Range('D47').Value = '3.90'
$ws.Range('E47').Value = '  -1.98%  '

# Row 48: update D48, E48
$ws.Range('D48').Value = '0.0586'
$ws.Range('E48').Value = '  -1.85%  '

# Row 49: update D49, E49
$ws.Range('D49').Value = '0.634'
$ws.Range('E49').Value = '  -0.25%  '

# Row 50: update B50, C50, D50, E50
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0248'
$ws.Range('E50').Value = '  -0.89%  '

# Row 51: update B51, C51, D51, E51
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.0996'
$ws.Range('E51').Value = '  -2.01%  '
